$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Original"
$ws.Name = "Original"

# Add the new URL value into G1 (new shared string)
$ws.Range("G1").Value = "https://www.linkedin.com/feed/update/urn:li:activity:7240826648583962625/"

# Update the active selection to J9 (was J21)
$ws.Range("J9").Select()
